$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.093421219042667758
$ws.Range("A2").Value = -0.0099999996805735236
$ws.Range("A3").Value = -0.044865883908713045
$ws.Range("A4").Value = 0.28399428578381247
$ws.Range("A5").Value = -0.0059999996928228327
$ws.Range("A6").Value = -0.0059999996811832546
$ws.Range("A7").Value = -0.019999999619765063
$ws.Range("A8").Value = -0.019999999616083564
$ws.Range("A9").Value = -0.0059999996720465631
$ws.Range("A10").Value = -0.0059999996684041434
$ws.Range("A11").Value = -0.0044999996748700255
$ws.Range("A12").Value = -0.005999999667515965
$ws.Range("A13").Value = -0.0059999996642021713
$ws.Range("A14").Value = 0.003772726859088138
$ws.Range("A15").Value = -0.0059999996620918594
$ws.Range("A16").Value = -0.0059999996609385597
$ws.Range("A17").Value = -0.0059999996594637395
$ws.Range("A18").Value = -0.0089999996462131193
$ws.Range("A19").Value = -0.0089999996846326091
$ws.Range("A20").Value = -0.0089999996818388439
$ws.Range("A21").Value = -0.0089999996814151828
$ws.Range("A22").Value = -0.008999999681112314
$ws.Range("A23").Value = -0.0089999996801548576
$ws.Range("A24").Value = -0.041999999533167021
$ws.Range("A25").Value = -0.041999999530694332
$ws.Range("A26").Value = -0.0059999996802240219
$ws.Range("A27").Value = -0.0059999996787505339
$ws.Range("A28").Value = -0.0059999996729862559
$ws.Range("A29").Value = -0.0026365076628813711
$ws.Range("A30").Value = -0.019999999606655106
$ws.Range("A31").Value = -0.014999999624970783
$ws.Range("A32").Value = -0.020999999598624086
$ws.Range("A33").Value = -0.0059999996629569452
